$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.225.02"

$ws.Range("D3").Value = "1.582.82"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").Value = "209.83"
$ws.Range("E5").Value = "  -0.75%  "

$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  -1.47%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "0.0611"
$ws.Range("E8").Value = "  -1.31%  "

$ws.Range("E9").Value = "  -0.24%  "

$ws.Range("D10").Value = "19.55"
$ws.Range("E10").Value = "  -0.44%  "

$ws.Range("E11").Value = "  +0.23%  "

$ws.Range("D12").Value = "1.805.91"
$ws.Range("E12").Value = "  -1.10%  "

$ws.Range("D13").Value = "1.577.60"
$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("E15").Value = "  -0.94%  "

$ws.Range("D16").Value = "'64.70"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("D17").Value = "26.221.73"
$ws.Range("E17").Value = "  -1.74%  "

$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").Value = "7.22"
$ws.Range("E19").Value = "  +0.40%  "

$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").Value = "206.99"
$ws.Range("E21").Value = "  -1.47%  "

$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("E23").Value = "  -3.29%  "

$ws.Range("D24").Value = "8.85"
$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("D25").Value = "144.65"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").Value = "7.01"
$ws.Range("E27").Value = "  -1.29%  "

$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("D29").Value = "15.25"
$ws.Range("E29").Value = "  -0.70%  "

$ws.Range("E30").Value = "  -1.84%  "

$ws.Range("E31").Value = "  -0.91%  "

$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("D33").Value = "2.94"
$ws.Range("E33").Value = "  -1.11%  "

$ws.Range("D34").Value = "1.28"
$ws.Range("E34").Value = "  +8.10%  "

$ws.Range("D35").Value = "1.288.00"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("E36").Value = "  -0.06%  "

$ws.Range("D37").Value = "0.604"
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("E38").Value = "  -0.95%  "

$ws.Range("D39").Value = "0.0166"
$ws.Range("E39").Value = "  -1.59%  "

$ws.Range("D40").Value = "0.815"
$ws.Range("E40").Value = "  -1.13%  "

$ws.Range("D41").Value = "5.53"
$ws.Range("E41").Value = "  +2.42%  "

$ws.Range("D42").Value = "0.769"
$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("D43").Value = "2.13"
$ws.Range("E43").Value = "  -2.70%  "

$ws.Range("D44").Value = "62.48"
$ws.Range("E44").Value = "  -0.88%  "

$ws.Range("D45").Value = "1.718.55"
$ws.Range("E45").Value = "  -1.20%  "

$ws.Range("D46").Value = "88.81"
$ws.Range("E46").Value = "  -2.00%  "

$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.42"
$ws.Range("E51").Value = "  -0.09%  "
